# AgileReleaseSummary.xlsx edit
# "penultimate doc commit still awaiting some analysis docs"
#
# Summary of the change:
#  - The old row 44 ("34 - Report Designs") is removed (its analysis doc is
#    done and merged away), shifting the trailing separator / totals rows
#    up by one.
#  - Three "banner" rows (12, 28 and the new 44) are restyled to match the
#    gray banner look already used by the totals row.
#  - A new "% Complete" tracking column (H) is populated for the
#    in-progress items (rows 29-43) with either an "X" marker or a
#    completion percentage.
#  - A couple of stray Work-Effort hours are filled in (E39, E40).
#  - The Work Effort totals row formulas are refreshed for the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Remove the old row 44 ("34  Report Designs"). Everything below it
#    (the two separator rows, the blank spacer row and the totals row)
#    shifts up by one, and Excel automatically re-points every SUM()
#    range that referenced rows below the deletion point.
# ---------------------------------------------------------------------
$ws.Rows.Item(44).Delete()

# ---------------------------------------------------------------------
# 2. Re-style the banner / separator rows (12, 28 and the new 44) so they
#    match the gray banner styling already used on the totals row, and
#    make sure any leftover text in them is cleared.
# ---------------------------------------------------------------------
$ws.Range("A47:G47").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("A28:G28").PasteSpecial(-4122)
$ws.Range("A44:G44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A12:G12").ClearContents()
$ws.Range("A28:G28").ClearContents()
$ws.Range("A44:G44").ClearContents()

# ---------------------------------------------------------------------
# 3. New "% Complete" column H for the in-progress work items.
# ---------------------------------------------------------------------
$ws.Range("H29").Value = "X"
$ws.Range("H30").Value = 0.5
$ws.Range("H30").NumberFormat = "0%"
$ws.Range("H31").Value = "X"
$ws.Range("H32").Value = "X"
$ws.Range("H33").Value = "X"
$ws.Range("H34").Value = "X"
$ws.Range("H35").Value = 0
$ws.Range("H35").NumberFormat = "0%"
$ws.Range("H36").Value = 0
$ws.Range("H36").NumberFormat = "0%"
$ws.Range("H37").Value = "X"
$ws.Range("H38").Value = "X"
$ws.Range("H39").Value = "X"
$ws.Range("H40").Value = "X"
$ws.Range("H41").Value = "X"
$ws.Range("H42").Value = "X"
$ws.Range("H43").Value = 0.5
$ws.Range("H43").NumberFormat = "0%"

# A couple of work-effort hours that were logged for items 19/20.
$ws.Range("E39").Value = 1.5
$ws.Range("E40").Value = 5

# ---------------------------------------------------------------------
# 4. Refresh the Work Effort totals row. Row 47 is now the totals row
#    after the deletion above; the E column total only needs to cover the
#    in-progress analysis items (rows 29-43), not the whole sheet.
# ---------------------------------------------------------------------
$ws.Range("E47").Formula = "=SUM(E29:E43)"

# ---------------------------------------------------------------------
# 5. View state: the active cell moves on to K35 while scrolled down to
#    keep row 13 in view.
# ---------------------------------------------------------------------
try {
    $excel.ActiveWindow.ScrollRow = 13
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("K35").Select()
